$d = $word.ActiveDocument

$replacements = @(
    @("298÷7=42, 4", "347÷6=57, 5"),
    @("632÷8=79, 0", "171÷7=24, 3"),
    @("942÷6=157, 0", "821÷4=205, 1"),
    @("854÷6=142, 2", "209÷3=69, 2"),
    @("332÷7=47, 3", "276÷3=92, 0"),
    @("316÷2=158, 0", "223÷7=31, 6"),
    @("908÷3=302, 2", "499÷7=71, 2"),
    @("387÷6=64, 3", "868÷4=217, 0"),
    @("870÷5=174, 0", "543÷5=108, 3"),
    @("668÷2=334, 0", "414÷3=138, 0"),
    @("306÷6=51, 0", "841÷3=280, 1"),
    @("133÷8=16, 5", "260÷7=37, 1"),
    @("941÷8=117, 5", "492÷7=70, 2"),
    @("751÷3=250, 1", "165÷7=23, 4"),
    @("489÷5=97, 4", "698÷6=116, 2"),
    @("843÷5=168, 3", "200÷5=40, 0"),
    @("653÷4=163, 1", "879÷4=219, 3"),
    @("180÷4=45, 0", "174÷2=87, 0"),
    @("920÷9=102, 2", "548÷3=182, 2"),
    @("757÷6=126, 1", "821÷2=410, 1"),
    @("263÷6=43, 5", "915÷3=305, 0"),
    @("783÷2=391, 1", "629÷9=69, 8"),
    @("616÷6=102, 4", "742÷9=82, 4"),
    @("888÷6=148, 0", "782÷7=111, 5"),
    @("560÷4=140, 0", "539÷5=107, 4")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
